# JobRequestReport.xlsx edit:
# - Update DateFrom (K2) and DateTo (L2) values
# - Move the active cell selection to K2
# - Add list-based data validation rules to several columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date range values (serial dates, keeps existing date number format/style)
$ws.Range("K2").Value = 42401
$ws.Range("L2").Value = 43434

# Add dropdown (list) data validations, in the same order as the source workbook
$ws.Range("D2:D8").Validation.Add(3, 1, 1, """IT,BPO,PST,SSS-Shared Services,SHILOH,GC-IT,DIGITAL""")
$ws.Range("E2:E4").Validation.Add(3, 1, 1, """APL Logistics,Arbor Health,Arke""")
$ws.Range("C2:C4").Validation.Add(3, 1, 1, """India,Australia,Canada""")
$ws.Range("F2:F4").Validation.Add(3, 1, 1, """All,Active,Cancelled""")
$ws.Range("G2:G4").Validation.Add(3, 1, 1, """.Net,Ab Initio,Admin""")
$ws.Range("H2").Validation.Add(3, 1, 1, """Confirmed,Pipeline""")
$ws.Range("I2").Validation.Add(3, 1, 1, """All,Internal,External""")
$ws.Range("J2:J4").Validation.Add(3, 1, 1, """Spreeth B,Abilash N Gatti,Agnish  Ghosh""")
$ws.Range("M2:M4").Validation.Add(3, 1, 1, """Created On,Needed By,Modified On""")

# Move the selection to K2 (matches the sheetView selection in the diff)
$ws.Range("K2").Select()
